$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow value edits
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (A80)
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-24 for illustrative purposes only and are subject to change."
$ws.Range("A80").Value = $newText

# Update Weight (D) and Percent Change (E) values for each holding row
$ws.Range("D2").Value = [double]"0.06207811792476992"
$ws.Range("E2").Value = [double]"0.01331419915490706"
$ws.Range("D3").Value = [double]"0.03774466842869519"
$ws.Range("E3").Value = [double]"0.0130842813791725"
$ws.Range("D4").Value = [double]"0.03177955715913908"
$ws.Range("E4").Value = [double]"0.02288208182077733"
$ws.Range("D5").Value = [double]"0.02939450249613725"
$ws.Range("E5").Value = [double]"0.01124487366053684"
$ws.Range("D6").Value = [double]"0.02703372259897427"
$ws.Range("E6").Value = [double]"0.02916574038960307"
$ws.Range("D7").Value = [double]"0.02555685055888524"
$ws.Range("E7").Value = [double]"0.00541005778925352"
$ws.Range("D8").Value = [double]"0.1943729207366834"
$ws.Range("D9").Value = [double]"0.02484636194320051"
$ws.Range("E9").Value = [double]"-0.002398221806270473"
$ws.Range("D10").Value = [double]"0.02276643087410839"
$ws.Range("E10").Value = [double]"-0.0002898550724637072"
$ws.Range("D11").Value = [double]"0.02206521223372782"
$ws.Range("E11").Value = [double]"0.009320905459387685"
$ws.Range("D12").Value = [double]"0.01896102745263907"
$ws.Range("E12").Value = [double]"0.01107888631090503"
$ws.Range("D13").Value = [double]"0.02031852892090889"
$ws.Range("E13").Value = [double]"0.005424528301886777"
$ws.Range("D14").Value = [double]"0.01709324526752049"
$ws.Range("E14").Value = [double]"0.01831012778943353"
$ws.Range("D15").Value = [double]"0.01644367924377113"
$ws.Range("E15").Value = [double]"0.01645359169867588"
$ws.Range("D16").Value = [double]"0.01496833910638067"
$ws.Range("E16").Value = [double]"0.02195910484107988"
$ws.Range("D17").Value = [double]"0.0144764805014711"
$ws.Range("E17").Value = [double]"7.054673721329507E-05"
$ws.Range("D18").Value = [double]"0.01444505685637226"
$ws.Range("E18").Value = [double]"0.005098572399728107"
$ws.Range("D19").Value = [double]"0.01366351152320778"
$ws.Range("E19").Value = [double]"0.02656294469215426"
$ws.Range("D20").Value = [double]"0.01296036818456489"
$ws.Range("E20").Value = [double]"0.01171079429735222"
$ws.Range("D21").Value = [double]"0.01166991691902473"
$ws.Range("E21").Value = [double]"-0.008663778740419903"
$ws.Range("D22").Value = [double]"0.01322731204967953"
$ws.Range("E22").Value = [double]"-0.001808473992421678"
$ws.Range("D23").Value = [double]"0.01155133193833298"
$ws.Range("E23").Value = [double]"0.01300326441784549"
$ws.Range("D24").Value = [double]"0.01293491503203483"
$ws.Range("E24").Value = [double]"-0.001651969001287612"
$ws.Range("D25").Value = [double]"0.01116297496446771"
$ws.Range("E25").Value = [double]"-0.002945181619533166"
$ws.Range("D26").Value = [double]"0.008877808213323814"
$ws.Range("E26").Value = [double]"0.02700693756194261"
$ws.Range("D27").Value = [double]"0.00964863022759832"
$ws.Range("E27").Value = [double]"0.03344732128317851"
$ws.Range("D28").Value = [double]"0.009917066715855146"
$ws.Range("E28").Value = [double]"0.01102687801516189"
$ws.Range("D29").Value = [double]"0.01010843671450707"
$ws.Range("E29").Value = [double]"-0.006077421059585086"
$ws.Range("D30").Value = [double]"0.009709906335541054"
$ws.Range("E30").Value = [double]"0.01626213592233006"
$ws.Range("D31").Value = [double]"0.008590910333571421"
$ws.Range("E31").Value = [double]"0.01569186875891582"
$ws.Range("D32").Value = [double]"0.009952182639253097"
$ws.Range("E32").Value = [double]"0.01053013798111824"
$ws.Range("D33").Value = [double]"0.009165138168196318"
$ws.Range("E33").Value = [double]"0.001230012300122985"
$ws.Range("D34").Value = [double]"0.009010887350317395"
$ws.Range("E34").Value = [double]"0.003295496155254485"
$ws.Range("D35").Value = [double]"0.00911246428309939"
$ws.Range("E35").Value = [double]"0.0003448424501055758"
$ws.Range("D36").Value = [double]"0.00828146598846061"
$ws.Range("E36").Value = [double]"0.001366003585759445"
$ws.Range("D37").Value = [double]"0.008450800155986976"
$ws.Range("E37").Value = [double]"0.001464128843338131"
$ws.Range("D38").Value = [double]"0.00684501261187996"
$ws.Range("E38").Value = [double]"0.0440022035532297"
$ws.Range("D39").Value = [double]"0.008787622351890149"
$ws.Range("E39").Value = [double]"-0.003504380475594537"
$ws.Range("D40").Value = [double]"0.007770596078266245"
$ws.Range("E40").Value = [double]"0.01421436803688025"
$ws.Range("D41").Value = [double]"0.007089802807199912"
$ws.Range("E41").Value = [double]"0.01761811896108489"
$ws.Range("D42").Value = [double]"0.00722916667321326"
$ws.Range("E42").Value = [double]"0.01102997109386883"
$ws.Range("D43").Value = [double]"0.008071084684523888"
$ws.Range("E43").Value = [double]"0.004428697962798989"
$ws.Range("D44").Value = [double]"0.00729515632792082"
$ws.Range("E44").Value = [double]"0.01055329413538386"
$ws.Range("D45").Value = [double]"0.007264204037498464"
$ws.Range("E45").Value = [double]"0.0002595492494701368"
$ws.Range("D46").Value = [double]"0.007903675215259826"
$ws.Range("E46").Value = [double]"-0.001192748091603191"
$ws.Range("D47").Value = [double]"0.007464372656778065"
$ws.Range("E47").Value = [double]"-0.001894417782268287"
$ws.Range("D48").Value = [double]"0.007135288533480481"
$ws.Range("E48").Value = [double]"0.0007596859964547065"
$ws.Range("D49").Value = [double]"0.006535332589430908"
$ws.Range("E49").Value = [double]"0.01442481067435986"
$ws.Range("D50").Value = [double]"0.007760933307398353"
$ws.Range("E50").Value = [double]"-0.001113461752588973"
$ws.Range("D51").Value = [double]"0.006523077367842362"
$ws.Range("E51").Value = [double]"0.007629404822121044"
$ws.Range("D52").Value = [double]"0.006762132747931774"
$ws.Range("E52").Value = [double]"0.01346468859277161"
$ws.Range("D53").Value = [double]"0.005327486230944309"
$ws.Range("E53").Value = [double]"0.02123424021234244"
$ws.Range("D54").Value = [double]"0.006280958182355812"
$ws.Range("E54").Value = [double]"0.005703422053231932"
$ws.Range("D55").Value = [double]"0.005333967357745945"
$ws.Range("E55").Value = [double]"0.006296255384955263"
$ws.Range("D56").Value = [double]"0.00568100223830624"
$ws.Range("E56").Value = [double]"0.005368180875337192"
$ws.Range("D57").Value = [double]"0.006657727687090885"
$ws.Range("E57").Value = [double]"0.007363005616651819"
$ws.Range("D58").Value = [double]"0.00556198518249439"
$ws.Range("E58").Value = [double]"-0.001694915254237372"
$ws.Range("D59").Value = [double]"0.005490103594330798"
$ws.Range("E59").Value = [double]"0.004893754024468677"
$ws.Range("D60").Value = [double]"0.004962422034008554"
$ws.Range("E60").Value = [double]"-0.007345491388044523"
$ws.Range("D61").Value = [double]"0.004774430077204754"
$ws.Range("E61").Value = [double]"0.009954751131221684"
$ws.Range("D62").Value = [double]"0.004903267022109991"
$ws.Range("E62").Value = [double]"0.01145557958824006"
$ws.Range("D63").Value = [double]"0.004158919428831257"
$ws.Range("E63").Value = [double]"-0.001737816395919856"
$ws.Range("D64").Value = [double]"0.004213910807754224"
$ws.Range("E64").Value = [double]"0.003579418344519247"
$ws.Range("D65").Value = [double]"0.003841226376882001"
$ws.Range("E65").Value = [double]"-0.003108638743455461"
$ws.Range("D66").Value = [double]"0.003803360884537901"
$ws.Range("E66").Value = [double]"-0.001611104226050264"
$ws.Range("D67").Value = [double]"0.003803360884537901"
$ws.Range("E67").Value = [double]"0.005019209319618412"
$ws.Range("D68").Value = [double]"0.003648285195975134"
$ws.Range("E68").Value = [double]"0.007235142118863003"
$ws.Range("D69").Value = [double]"0.003503736428520478"
$ws.Range("E69").Value = [double]"-0.007197309417040332"
$ws.Range("D70").Value = [double]"0.00292774101385877"
$ws.Range("E70").Value = [double]"0.004132231404958775"
$ws.Range("D71").Value = [double]"0.002954097596185421"
$ws.Range("E71").Value = [double]"0.02584865770473499"
$ws.Range("D72").Value = [double]"0.00235547715705255"
$ws.Range("E72").Value = [double]"0.0413727550152585"
$ws.Range("D73").Value = [double]"0.001955689832282579"
$ws.Range("E73").Value = [double]"0.01006246359637664"
$ws.Range("D74").Value = [double]"0.001917117307923755"
$ws.Range("E74").Value = [double]"0.01999713155899774"
$ws.Range("D75").Value = [double]"0.001494115765337017"
$ws.Range("E75").Value = [double]"0.02686786897313209"
$ws.Range("D76").Value = [double]"0.001671423682807212"
$ws.Range("E76").Value = [double]"-0.0002820078962211259"
$ws.Range("D77").Value = [double]"1"
$ws.Range("E77").Value = [double]"0.009713692884775416"
